$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'29.167.95"
$ws.Range("E2").Value2 = "'  -2.91%  "

$ws.Range("D3").Value2 = "'1.848.87"
$ws.Range("E3").Value2 = "'  -1.90%  "

$ws.Range("D4").Value2 = "'1.000"
$ws.Range("E4").Value2 = "'  -0.29%  "

$ws.Range("D5").Value2 = "'0.7032"
$ws.Range("E5").Value2 = "'  -4.57%  "

$ws.Range("D6").Value2 = "'238.64"
$ws.Range("E6").Value2 = "'  -1.69%  "

$ws.Range("E7").Value2 = "'  -0.21%  "

$ws.Range("E8").Value2 = "'  -3.52%  "

$ws.Range("D9").Value2 = "'0.07473"
$ws.Range("E9").Value2 = "'  +4.16%  "

$ws.Range("D10").Value2 = "'23.45"
$ws.Range("E10").Value2 = "'  -4.88%  "

$ws.Range("D11").Value2 = "'0.08130"
$ws.Range("E11").Value2 = "'  -2.61%  "

$ws.Range("B12").Value2 = "'Polygon"
$ws.Range("C12").Value2 = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value2 = "'0.7268"
$ws.Range("E12").Value2 = "'  -3.84%  "

$ws.Range("B13").Value2 = "'Polkadot"
$ws.Range("C13").Value2 = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value2 = "'5.233"
$ws.Range("E13").Value2 = "'  -3.13%  "

$ws.Range("B14").Value2 = "'WrappedEther"
$ws.Range("C14").Value2 = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value2 = "'1.832.22"
$ws.Range("E14").Value2 = "'  -4.15%  "

$ws.Range("D15").Value2 = "'88.81"
$ws.Range("E15").Value2 = "'  -4.25%  "

$ws.Range("D16").Value2 = "'29.133.60"
$ws.Range("E16").Value2 = "'  -2.81%  "

$ws.Range("D17").Value2 = "'5.773"
$ws.Range("E17").Value2 = "'  -6.00%  "

$ws.Range("D18").Value2 = "'238.37"
$ws.Range("E18").Value2 = "'  -4.45%  "

$ws.Range("E19").Value2 = "'  -3.65%  "

$ws.Range("D20").Value2 = "'0.000007623"
$ws.Range("E20").Value2 = "'  -2.94%  "

$ws.Range("D21").Value2 = "'0.9990"
$ws.Range("E21").Value2 = "'  -0.06%  "

$ws.Range("D22").Value2 = "'2.092.32"
$ws.Range("E22").Value2 = "'  -2.95%  "

$ws.Range("D23").Value2 = "'1.000"
$ws.Range("E23").Value2 = "'  -0.25%  "

$ws.Range("D24").Value2 = "'7.593"
$ws.Range("E24").Value2 = "'  -3.59%  "

$ws.Range("D25").Value2 = "'9.004"
$ws.Range("E25").Value2 = "'  -2.96%  "

$ws.Range("D26").Value2 = "'161.27"
$ws.Range("E26").Value2 = "'  -1.86%  "

$ws.Range("D27").Value2 = "'0.1450"
$ws.Range("E27").Value2 = "'  -7.25%  "

$ws.Range("D28").Value2 = "'18.07"
$ws.Range("E28").Value2 = "'  -3.02%  "

$ws.Range("D29").Value2 = "'1.979"
$ws.Range("E29").Value2 = "'  -3.21%  "

$ws.Range("D30").Value2 = "'1.396"
$ws.Range("E30").Value2 = "'  -5.46%  "

$ws.Range("E31").Value2 = "'  +0.08%  "

$ws.Range("D32").Value2 = "'1.495"
$ws.Range("E32").Value2 = "'  -2.71%  "

$ws.Range("D33").Value2 = "'3.984"
$ws.Range("E33").Value2 = "'  -4.96%  "

$ws.Range("E34").Value2 = "'  -2.93%  "

$ws.Range("E35").Value2 = "'  -4.95%  "

$ws.Range("D36").Value2 = "'1.032"
$ws.Range("E36").Value2 = "'  +3.57%  "

$ws.Range("D37").Value2 = "'0.7036"
$ws.Range("E37").Value2 = "'  -8.17%  "

$ws.Range("D38").Value2 = "'2.656"
$ws.Range("E38").Value2 = "'  -2.73%  "

$ws.Range("D39").Value2 = "'0.01867"
$ws.Range("E39").Value2 = "'  -4.59%  "

$ws.Range("D40").Value2 = "'2.680"
$ws.Range("E40").Value2 = "'  -2.96%  "

$ws.Range("D41").Value2 = "'0.9394"
$ws.Range("E41").Value2 = "'  +7.51%  "

$ws.Range("D42").Value2 = "'6.016"
$ws.Range("E42").Value2 = "'  -0.52%  "

$ws.Range("D43").Value2 = "'1.074.97"
$ws.Range("E43").Value2 = "'  -2.39%  "

$ws.Range("D44").Value2 = "'0.4288"
$ws.Range("E44").Value2 = "'  -5.76%  "

$ws.Range("D45").Value2 = "'70.23"
$ws.Range("E45").Value2 = "'  -2.78%  "

$ws.Range("D46").Value2 = "'0.9997"
$ws.Range("E46").Value2 = "'  -0.41%  "

$ws.Range("D47").Value2 = "'102.69"
$ws.Range("E47").Value2 = "'  -1.35%  "

$ws.Range("D48").Value2 = "'1.741"
$ws.Range("E48").Value2 = "'  -6.05%  "

$ws.Range("D49").Value2 = "'1.988.53"
$ws.Range("E49").Value2 = "'  -2.53%  "

$ws.Range("B50").Value2 = "'Aptos"
$ws.Range("C50").Value2 = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value2 = "'7.055"
$ws.Range("E50").Value2 = "'  -6.74%  "

$ws.Range("B51").Value2 = "'EnergySwap"
$ws.Range("C51").Value2 = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value2 = "'9.160"
$ws.Range("E51").Value2 = "'  -4.14%  "

